# Revert "changed decimal separator from comma to dot"
#
# The egg-definitions table (sheet "gacha", table `eggDefinitions`) has
# three weight columns - [weightCommon] (F), [weightRare] (G) and
# [weightEpic] (H) - for rows 5-9. These cells were left blank by the
# previous (reverted) edit; restore their original numeric weights of
# 1, 2 and 3 respectively, matching rows 10-12 which already carry
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gacha")

$rows = 5,6,7,8,9
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 1   # F -> [weightCommon]
    $ws.Cells.Item($r, 7).Value = 2   # G -> [weightRare]
    $ws.Cells.Item($r, 8).Value = 3   # H -> [weightEpic]
}

# Restore the on-screen selection/scroll state left by the author after
# making the edit (row 20, topLeftCell A19 in the original commit).
$ws.Activate()
$ws.Range("J20").Select()
